# Applies the "Handles float input without breaking stuff" edit to the
# marksheet worksheet: updates the summary numbers (rows 10-12), removes the
# unused 2nd/3rd "Student Ans / Correct Ans" columns beyond row 18, and marks
# the student's recorded answer (column A / D) equal to the correct answer
# wherever the diff shows a match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# --- Summary block (rows 10-12) ------------------------------------------

# Give the row-label cells (A10/A11/A12) the same bold boxed style already
# used by the other header labels (copy format only, keep existing text).
# NOTE: PasteSpecial only honours the first area of a multi-area (union)
# range in this runtime, so paste into each target cell individually.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B10").Value = 21
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 84
$ws.Range("E12").Value = "84/112"

# --- Drop the 2nd/3rd "Student Ans / Correct Ans" blocks ------------------

# Columns F:H are only ever used by the 3rd block (rows 15-21); clearing the
# whole rectangle removes it completely (and shrinks the sheet dimension).
$ws.Range("F5:H40").Clear()

# The 2nd block (columns D:E) is kept only for rows 16-18; clear the rest.
$ws.Range("D19:E40").Clear()

# --- Record the student's answer where it equals the correct answer -------

# Re-use the existing "correct" style (green, centered, boxed) already
# applied to B10/B11/B12 instead of letting Excel synthesize a new one.
# NOTE: PasteSpecial only honours the first area of a multi-area (union)
# range in this runtime, so paste into each target cell individually.
$correctStyleCells = "A16","A18","A19","A20","A21","A22","A27","A28","A29","A30","A31","A32","A33","A35","A36","A37","A38","A39","D16","D17","D18"
foreach ($cellRef in $correctStyleCells) {
    $ws.Range("B10").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "Option A"
$ws.Range("A18").Value = "Option B"
$ws.Range("A19").Value = "Option C"
$ws.Range("A20").Value = "Option B"
$ws.Range("A21").Value = "Option C"
$ws.Range("A22").Value = "Option D"
$ws.Range("A27").Value = "Option A"
$ws.Range("A28").Value = "Option D"
$ws.Range("A29").Value = "Option D"
$ws.Range("A30").Value = "Option B"
$ws.Range("A31").Value = "Option D"
$ws.Range("A32").Value = "Option C"
$ws.Range("A33").Value = "Option D"
$ws.Range("A35").Value = "Option D"
$ws.Range("A36").Value = "Option A"
$ws.Range("A37").Value = "Option A"
$ws.Range("A38").Value = "Option A"
$ws.Range("A39").Value = "Option D"

$ws.Range("D16").Value = "Option A"
$ws.Range("D17").Value = "Option C"
$ws.Range("D18").Value = "Option D"
